$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 179; existing rows 179..263 shift down to 180..264.
$ws.Rows.Item(179).Insert()

# Populate the newly inserted row 179 with the new weekly data point.
$ws.Cells.Item(179, 1).Value = 3
$ws.Cells.Item(179, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(179, 3).Value = "Coquimbo"
$ws.Cells.Item(179, 4).Value = 45229
$ws.Cells.Item(179, 5).Value = 5
$ws.Cells.Item(179, 6).Value = 100112052
$ws.Cells.Item(179, 7).Value = "Albahaca"
$ws.Cells.Item(179, 8).Value = "Sin especificar"
$ws.Cells.Item(179, 9).Value = "Primera"
$ws.Cells.Item(179, 10).Value = 60
$ws.Cells.Item(179, 11).Value = 5000
$ws.Cells.Item(179, 12).Value = 5000
$ws.Cells.Item(179, 13).Value = 5000
$ws.Cells.Item(179, 14).Value = "$/docena de matas"
$ws.Cells.Item(179, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(179, 16).Value = 833
$ws.Cells.Item(179, 17).Value = 6
$ws.Cells.Item(179, 18).Value = "Hortaliza"
